# studentData.xlsx - "Added testing to project"
#
# Appends one more EmailKey / Fullname row to the student roster sheet:
#   A23 = s@nwmissouri.edu   (becomes a mailto: hyperlink, like the other
#                              EmailKey cells would if Excel auto-linked them)
#   B23 = Sravya Kancharla   (re-uses the existing "Sravya Kancharla" shared string)
#
# Adding the hyperlink pulls in Excel's built-in "Hyperlink" cell style
# (an underlined, theme-colored font) the first time it's used in the
# workbook, which is why styles.xml gains a second font/xf/cellStyle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "s@nwmissouri.edu"
$ws.Range("B23").Value = "Sravya Kancharla"

$ws.Hyperlinks.Add($ws.Range("A23"), "mailto:s@nwmissouri.edu") | Out-Null

# Leave the cursor where the author's saved view shows it.
$ws.Range("A8").Select() | Out-Null
